$d = $word.ActiveDocument

# 1) Merge the two runs "E" and "-NSW" (around the "Manhattan phase id: 6" text)
#    into a single run "E-NSW". Doing the replace across the run boundary via
#    Find/Replace both merges the text into one run and removes the
#    bookmarkStart/bookmarkEnd ("_GoBack") that used to sit between them.
$null = $d.Content.Find.Execute("E-NSW", $false, $false, $false, $false, $false, `
                                 $true, 1, $false, "E-NSW", 2)

# 2) Re-insert the "_GoBack" bookmark so that it wraps just the "N" run that
#    follows "Manhattan phase id: 8: " (i.e. before "N-WES").
$r = $d.Content
$null = $r.Find.Execute("Manhattan phase id: 8: ", $false, $false, $false, `
                         $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$null = $r.Find.Execute("N", $false, $false, $false, $false, $false, $true, `
                         1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $r)
